$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# paises.xlsx / "Pais" sheet update ("Update countries & provincias Spain")
#
# The country list (sorted descending by "Casos totales", col B) is
# re-ranked because Senegal's numbers jumped above "Republica de Yibuti"
# and "Hong Kong". That pushes those two countries down one row each
# (rows 89-91), while their own totals are unchanged. A handful of other
# countries simply received refreshed daily figures (no re-sort needed).
# ---------------------------------------------------------------------------

# --- Sanity check current layout before touching rows 89-91 ---
$a89 = $ws.Range("A89").Value()
$a90 = $ws.Range("A90").Value()
$a91 = $ws.Range("A91").Value()
if (($a89 -eq "Republica de Yibuti") -and ($a90 -eq "Hong Kong") -and ($a91 -eq "Senegal")) {
    Write-Host "Pre-check OK: rows 89-91 are Yibuti/HongKong/Senegal as expected"
} else {
    Write-Host "Pre-check WARNING: unexpected rows 89-91 -> $a89 / $a90 / $a91"
}

# --- Re-rank rows 89-91: Senegal moves up, Yibuti and Hong Kong move down ---
$ws.Range("A89").Value = "Senegal"
$ws.Range("B89").Value = 1115
$ws.Range("C89").Value = 91
$ws.Range("D89").Value = 368
$ws.Range("E89").Value = 738
$ws.Range("F89").Value = 1
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 9

$ws.Range("A90").Value = "Republica de Yibuti"
$ws.Range("B90").Value = 1112
$ws.Range("C90").Value = 15
$ws.Range("D90").Value = 686
$ws.Range("E90").Value = 424
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 2

$ws.Range("A91").Value = "Hong Kong"
$ws.Range("B91").Value = 1040
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 864
$ws.Range("E91").Value = 172
$ws.Range("F91").Value = 4
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 4

# --- Refresh daily figures for the remaining changed countries ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1135657
$ws.Range("C4").Value = 4627
$ws.Range("E4").Value = 907932
$ws.Range("G4").Value = 190
$ws.Range("H4").Value = 65943

# Row 9 - Alemania
$ws.Range("B9").Value = 164316
$ws.Range("C9").Value = 239
$ws.Range("E9").Value = 28580

# Row 15 - Canada
$ws.Range("B15").Value = 55572
$ws.Range("C15").Value = 511
$ws.Range("D15").Value = 23316
$ws.Range("E15").Value = 28810
$ws.Range("G15").Value = 55
$ws.Range("H15").Value = 3446

# Row 20 - Suiza
$ws.Range("E20").Value = 4157
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 1760

# Row 42 - Serbia
$ws.Range("F42").Value = 57

# Row 58 - Argelia
$ws.Range("B58").Value = 4295
$ws.Range("C58").Value = 141
$ws.Range("D58").Value = 1872
$ws.Range("E58").Value = 1964
$ws.Range("G58").Value = 6
$ws.Range("H58").Value = 459

# Row 65 - Grecia
$ws.Range("B65").Value = 2620
$ws.Range("C65").Value = 8
$ws.Range("E65").Value = 1103
$ws.Range("F65").Value = 37
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = 143

# Row 72 - Uzbekistan
$ws.Range("B72").Value = 2113
$ws.Range("C72").Value = 27
$ws.Range("E72").Value = 833

# Row 95 - Republica de Chipre
$ws.Range("B95").Value = 864
$ws.Range("C95").Value = 7
$ws.Range("E95").Value = 553

# Row 126 - Mauricio
$ws.Range("D126").Value = 314
$ws.Range("E126").Value = 8

# Row 150 - Trinidad yTobago
$ws.Range("D150").Value = 87
$ws.Range("E150").Value = 21

Write-Host "paises.xlsx update complete"
